$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New frequency values for B2:X5 (one row per data row, in column order B..X)
$data = @(
    @(0.00666173205033309, 0.0111028867505551, 0.00296076980014804, 0.00148038490007402, 0.0111028867505551, 0.117690599555885, 0.0399703923019985, 0.0162842339008142, 0.00370096225018505, 0.0125832716506292, 0.00518134715025907, 0.00370096225018505, 0.00370096225018505, 0.0185048112509252, 0.00148038490007402, 0.00666173205033309, 0.977054034048853, 0.0074019245003701, 0.0429311621021466, 0.0266469282013323, 0.00814211695040711, 0.00148038490007402, 0.0125832716506292),
    @(0.00222057735011103, 0.686158401184308, 0.00222057735011103, 0.0421909696521095, 0.00074019245003701, 0.00888230940044412, 0.920059215396003, 0.973353071798668, 0.00074019245003701, 0.00074019245003701, 0.00666173205033309, 0.00592153960029608, 0.974833456698742, 0.00074019245003701, 0, 0, 0.00148038490007402, 0.982975573649149, 0.0895632864544782, 0.00518134715025907, 0, 0.00962250185048113, 0.00370096225018505),
    @(0.990377498149519, 0.00370096225018505, 0.988156920799408, 0.0133234641006662, 0.982235381199112, 0.86380458919319, 0.00296076980014804, 0.00074019245003701, 0.0259067357512953, 0.981495188749075, 0.987416728349371, 0.988156920799408, 0.00222057735011103, 0.978534418948927, 0.998519615099926, 0.99259807549963, 0.0162842339008142, 0.00074019245003701, 0.0259067357512953, 0.962250185048113, 0.983715766099186, 0.987416728349371, 0.97779422649889),
    @(0.00074019245003701, 0.299037749814952, 0.00666173205033309, 0.94300518134715, 0.00592153960029608, 0.00888230940044412, 0.0362694300518135, 0.00962250185048113, 0.969652109548483, 0.00444115470022206, 0.00074019245003701, 0.00074019245003701, 0.0192450037009623, 0.00148038490007402, 0, 0.00074019245003701, 0.00518134715025907, 0.00888230940044412, 0.84159881569208, 0.00444115470022206, 0, 0.00148038490007402, 0.00444115470022206)
)

$firstRow = 2
$firstCol = 2   # column B

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($firstRow + $r, $firstCol + $c).Value = $rowValues[$c]
    }
}
